$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.089.26'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '4.023.96'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '531.70'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.96'
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.700'
$ws.Range('E7').Value = '  +11.56%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.749'
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('E10').Value = '  -3.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000326'
$ws.Range('E11').Value = '  -5.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.76'
$ws.Range('E12').Value = '  +1.89%  '
$ws.Range('D13').Value = '4.685.02'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.63'
$ws.Range('E14').Value = '  -1.78%  '
$ws.Range('D15').Value = '4.029.39'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.11'
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.54'
$ws.Range('E17').Value = '  -4.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.132'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('E19').Value = '  -3.16%  '
$ws.Range('D20').Value = '72.071.92'
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '429.36'
$ws.Range('E21').Value = '  -3.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '97.72'
$ws.Range('E22').Value = '  +2.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.46'
$ws.Range('E23').Value = '  -3.40%  '
$ws.Range('E24').Value = '  +3.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.11'
$ws.Range('E25').Value = '  -1.97%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.09'
$ws.Range('E26').Value = '  -10.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.69'
$ws.Range('E27').Value = '  -4.46%  '
$ws.Range('E28').Value = '  +0.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.71'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.58'
$ws.Range('E30').Value = '  +22.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.32'
$ws.Range('E31').Value = '  -2.55%  '
$ws.Range('E32').Value = '  -2.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '675.40'
$ws.Range('E33').Value = '  -3.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.06'
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '44.45'
$ws.Range('E35').Value = '  +8.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '65.93'
$ws.Range('E36').Value = '  -2.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.440'
$ws.Range('E37').Value = '  -1.00%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0830'
$ws.Range('E38').Value = '  -9.95%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.152'
$ws.Range('E39').Value = '  -3.17%  '
$ws.Range('E40').Value = '  -4.47%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0484'
$ws.Range('E43').Value = '  -1.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.18'
$ws.Range('E44').Value = '  +1.89%  '
$ws.Range('E45').Value = '  +2.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.72'
$ws.Range('E46').Value = '  +5.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.43'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.61'
$ws.Range('E48').Value = '  -7.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.00'
$ws.Range('E49').Value = '  -7.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000269'
$ws.Range('E50').Value = '  -3.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.59'
$ws.Range('E51').Value = '  +0.36%  '
